# Updated legacy GSC export data.
# The first data row (2025-11-04) is removed, the remaining daily rows
# shift up by one, and the newly-introduced first three rows
# (2025-11-05 .. 2025-11-07) lose their "Not indexed"/"Indexed" values
# (now blank) while keeping the "Impressions" carried from the old data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Delete the old first data row (row 2, date 2025-11-04). Excel shifts
# rows 3:90 up to 2:89 automatically, same as removing a row from a
# GSC export refresh.
$ws.Rows.Item(2).Delete()

# The newly top three data rows (now holding 2025-11-05, 2025-11-06,
# 2025-11-07) no longer carry "Not indexed"/"Indexed" figures for those
# dates in the refreshed export - clear them to blank.
$ws.Range("B2:C4").Value = ""
